$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing text values (shared-string edits) ---

# C4 (row for O3 "Registrare una spedizione"): "1 ogni 2 settimane" -> "1 ogni 10 giorni"
$ws.Range("C4").Value = "1 ogni 10 giorni"

# C5 (row for O4 "Aggiungere un avvistamento..."): "10 ogni 2 settimane" -> "10 ogni 10 giorni"
$ws.Range("C5").Value = "10 ogni 10 giorni"

# B10: shorten description text
$ws.Range("B10").Value = "Visualizzare i luoghi più pericolosi"

# Row 10 is no longer wrapped across 3 lines after the text shortened
$ws.Rows.Item(10).RowHeight = 18.7

# --- Add new row 13 (O12) ---

# Copy formatting from row 12 so the new row matches existing styles
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "O12"
$ws.Range("B13").Value = "Visualizzare quanti organismi vengono scoperti ogni anno"
$ws.Range("C13").Value = "1 all’anno"

$ws.Rows.Item(13).RowHeight = 31.3

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("F13").Select() | Out-Null
